# Add the new worksheet for device THA0821ASW01 as the last sheet in the workbook
$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "THA0821ASW01"

# Port table data (Port, Active, Name, PVID, FlowControl, Type, SpeedDuplex, Priority)
$data = @(
    @("Port", "Active", "Name", "PVID", "FlowControl", "Type", "SpeedDuplex", "Priority"),
    @("'1", "No", "'", "'1", "No", "100M", "100-full", "'0"),
    @("'2", "No", "'", "'1", "No", "100M", "100-full", "'0"),
    @("'3", "No", "'", "'1", "No", "100M", "100-full", "'0"),
    @("'4", "Yes", "t037_a_line_dinhbt1", "'1", "No", "100M", "100-full", "'0"),
    @("'5", "No", "'", "'1", "No", "100M", "100-full", "'0"),
    @("'6", "No", "'", "'1", "No", "100M", "100-full", "'0"),
    @("'7", "No", "'", "'1", "No", "100M", "100-full", "'0"),
    @("'8", "Yes", "t037_a_line_tuha", "'1", "No", "100M", "100-full", "'0"),
    @("'9", "No", "'", "'1", "No", "100M", "100-full", "'0"),
    @("'10", "No", "'", "'1", "No", "100M", "100-full", "'0"),
    @("'11", "No", "'", "'1", "No", "100M", "100-full", "'0"),
    @("'12", "No", "'", "'1", "No", "100M", "100-full", "'0"),
    @("'13", "No", "'", "'1", "No", "100M", "100-full", "'0"),
    @("'14", "No", "'", "'1", "No", "100M", "100-full", "'0"),
    @("'15", "No", "'", "'1", "No", "100M", "100-full", "'0"),
    @("'16", "No", "'", "'1", "No", "100M", "100-full", "'0"),
    @("'17", "No", "'", "'1", "No", "100M", "100-full", "'0"),
    @("'18", "No", "'", "'1", "No", "100M", "100-full", "'0"),
    @("'19", "No", "'", "'1", "No", "100M", "100-full", "'0"),
    @("'20", "No", "'", "'1", "No", "100M", "100-full", "'0"),
    @("'21", "No", "'", "'1", "No", "100M", "100-full", "'0"),
    @("'22", "No", "'", "'1", "No", "100M", "100-full", "'0"),
    @("'23", "No", "'", "'1", "No", "100M", "100-full", "'0"),
    @("'24", "No", "'", "'1", "No", "100M", "100-full", "'0"),
    @("'25", "Yes", "'", "'1", "No", "10/100/1000M", "auto", "'0"),
    @("'26", "Yes", "'", "'1", "No", "10/100/1000M", "auto", "'0"),
    @("'27", "Yes", "'", "'1", "No", "10/100/1000M", "auto", "'0"),
    @("'28", "Yes", "UPLINK_THA0821SRT03", "'1", "No", "10/100/1000M", "1000-full", "'0"),
)

for ($r = 0; $r -lt $data.Count; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}
